$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Jimmy Butler 2024-25 GSW) statistics totals
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 320
$ws.Range("E2").Value = 180
$ws.Range("F2").Value = 59
$ws.Range("G2").Value = 132
$ws.Range("H2").Value = 44.7
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 32
$ws.Range("K2").Value = 34.4
$ws.Range("L2").Value = 51
$ws.Range("M2").Value = 64
$ws.Range("N2").Value = 79.7
$ws.Range("O2").Value = 20
$ws.Range("P2").Value = 41
$ws.Range("Q2").Value = 61
$ws.Range("R2").Value = 48
$ws.Range("S2").Value = 10
$ws.Range("T2").Value = 11
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 15
$ws.Range("W2").Value = 354.2
$ws.Range("X2").Value = 1
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = -6

# Update the active cell selection on the sheet view
$ws.Range("X15").Select()
